# Ziegler-Nichols in-progress update: refresh simulation results (Kp, Tau, Td)
# for the ethanol-concentration cooling-flow step table, and re-write the header
# labels (Step, Kp, Tau, Td) whose shared-string entries were regenerated upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Step"
$ws.Range("B1").Value = "Kp"
$ws.Range("C1").Value = "Tau"
$ws.Range("D1").Value = "Td"

$data = @(
    @(-100, 0.0010109467996918652, 17.5805668441688, 8.0132943231523086),
    @(-95, 0.0009963862010532933, 17.584350455213553, 8.0161704430457235),
    @(-90, 0.00098222393252132487, 17.553194973368534, 8.0378466738270049),
    @(-85, 0.00096844427101613667, 17.595073160567779, 8.0348554740620557),
    @(-80, 0.0009550319506555383, 17.553701316752893, 8.0566369638042374),
    @(-75, 0.00094197268514686048, 17.559387772530954, 8.0610422603169809),
    @(-70, 0.00092925279013148356, 17.575331816466701, 8.0701838276324338),
    @(-65, 0.00091686042809846838, 17.632411728314935, 8.055796720089802),
    @(-60, 0.00090477939057507226, 17.645213561356911, 8.055590156756125),
    @(-55, 0.00089300129386112362, 17.67682875177195, 8.0537302667390804),
    @(-50, 0.00088151322788164064, 17.581329418597136, 8.0989366491597821),
    @(-45, 0.00087030389290971601, 17.710321743164343, 8.0603206440501936),
    @(-40, 0.00085936186428163812, 17.660260090397514, 8.0923605411482527),
    @(-35, 0.00084867560012733228, 17.721513468807302, 8.0727378824440734),
    @(-30, 0.00083823258655092063, 17.682021148999524, 8.0929530951611923),
    @(-25, 0.00082801812597750764, 17.756922962635258, 8.0765445778953904),
    @(-20, 0.00081801269223156223, 17.696780890335219, 8.109066513994037),
    @(-15, 0.00080817852404682828, 17.64223680541005, 8.134384766982123),
    @(-10, 0.00079842406530801209, 17.753238624472601, 8.110228747546671),
    @(10, 0.00076427604831224021, 17.823514576182617, 8.0891025535912036),
    @(15, 0.00075552196438195547, 17.825140082345399, 8.099106519770146),
    @(20, 0.0007470891668685198, 17.738867741327056, 8.130709299434443),
    @(25, 0.00073888799847225076, 17.86582021892724, 8.0958193118188433),
    @(30, 0.00073088498389755805, 17.85688167044998, 8.1062116162523807),
    @(35, 0.00072306292364276352, 17.761479833248586, 8.1430018622625084),
    @(40, 0.0007154114001474809, 17.79594866839966, 8.1406305462973023),
    @(45, 0.00070792061552853939, 17.875197000535422, 8.1232692054017548),
    @(50, 0.00070058560346282657, 17.875146629017706, 8.1233510798820134),
    @(55, 0.00069340012142605058, 17.82374567098492, 8.1454193737517926),
    @(60, 0.00068635906352116174, 17.787247508967937, 8.163277667301486),
    @(65, 0.00067945775294344474, 17.871983287053183, 8.1403592922997063),
    @(70, 0.00067269181290663079, 17.899933975662407, 8.137928902728504),
    @(75, 0.00066605713656249087, 17.816987300833546, 8.1754092781949907),
    @(80, 0.00065954986109391096, 17.925360059546392, 8.1416881484920509),
    @(85, 0.00065316625842797298, 17.89385534763943, 8.1557937387468655),
    @(90, 0.00064690281963572082, 17.932191316098852, 8.1534769225758055),
    @(95, 0.00064075615071257483, 17.808140824196471, 8.1989088999663977),
    @(100, 0.00063472302042754335, 17.897000012663796, 8.1725374942842564)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
